$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the header text in A1 and remove the bold/border style from row 1
$ws.Range("A1").Value = ""
$ws.Range("A1:BA1").ClearFormats()

# Update corrected numeric values (columns affected by the data-cleaning fix)
$ws.Range("D3").Value = 0
$ws.Range("F3").Value = 27
$ws.Range("G3").Value = 42
$ws.Range("H3").Value = 40
$ws.Range("J3").Value = 21
$ws.Range("X3").Value = 56
$ws.Range("Y3").Value = 38
$ws.Range("Z3").Value = 25
$ws.Range("AA3").Value = 38
$ws.Range("AB3").Value = 15
$ws.Range("AJ3").Value = 36
$ws.Range("AU3").Value = 6
$ws.Range("AW3").Value = 8
$ws.Range("D4").Value = 1
$ws.Range("F4").Value = 58
$ws.Range("G4").Value = 406
$ws.Range("H4").Value = 129
$ws.Range("J4").Value = 40
$ws.Range("X4").Value = 328
$ws.Range("Y4").Value = 83
$ws.Range("Z4").Value = 37
$ws.Range("AA4").Value = 67
$ws.Range("AB4").Value = 20
$ws.Range("AJ4").Value = 293
$ws.Range("AU4").Value = 10
$ws.Range("AW4").Value = 9
$ws.Range("D5").Value = 133.49
$ws.Range("F5").Value = 19696.89
$ws.Range("G5").Value = 123038.1
$ws.Range("H5").Value = 44810.45
$ws.Range("J5").Value = 16916.14
$ws.Range("X5").Value = 102627.1
$ws.Range("Y5").Value = 25142.93
$ws.Range("Z5").Value = 13588.07
$ws.Range("AA5").Value = 25756.74
$ws.Range("AB5").Value = 7366.75
$ws.Range("AJ5").Value = 123981.4
$ws.Range("AU5").Value = 4037.85
$ws.Range("AW5").Value = 3979
$ws.Range("B6").Value = 0.54
$ws.Range("C6").Value = 0.27
$ws.Range("D6").Value = 0.03
$ws.Range("E6").Value = 1.33
$ws.Range("F6").Value = 4.62
$ws.Range("G6").Value = 28.85
$ws.Range("H6").Value = 10.51
$ws.Range("I6").Value = 6.09
$ws.Range("J6").Value = 3.97
$ws.Range("K6").Value = 0.36
$ws.Range("L6").Value = 1.83
$ws.Range("M6").Value = 2.02
$ws.Range("P6").Value = 3.23
$ws.Range("Q6").Value = 2.08
$ws.Range("S6").Value = 0.85
$ws.Range("V6").Value = 0.13
$ws.Range("W6").Value = 1.28
$ws.Range("X6").Value = 24.06
$ws.Range("Y6").Value = 5.9
$ws.Range("Z6").Value = 3.19
$ws.Range("AA6").Value = 6.04
$ws.Range("AB6").Value = 1.73
$ws.Range("AC6").Value = 2.64
$ws.Range("AD6").Value = 2.5
$ws.Range("AE6").Value = 0.32
$ws.Range("AF6").Value = 0.45
$ws.Range("AH6").Value = 2.36
$ws.Range("AI6").Value = 1.15
$ws.Range("AJ6").Value = 29.07
$ws.Range("AK6").Value = 0.65
$ws.Range("AL6").Value = 1.3
$ws.Range("AM6").Value = 0.53
$ws.Range("AN6").Value = 0.2
$ws.Range("AO6").Value = 0.83
$ws.Range("AP6").Value = 1.11
$ws.Range("AR6").Value = 0.15
$ws.Range("AT6").Value = 0.41
$ws.Range("AU6").Value = 0.95
$ws.Range("AV6").Value = 2.12
$ws.Range("AW6").Value = 0.93
$ws.Range("AX6").Value = 0.79
$ws.Range("AZ6").Value = 1.13
$ws.Range("BA6").Value = 1.62
$ws.Range("D7").Value = 133.49
$ws.Range("F7").Value = 339.6
$ws.Range("G7").Value = 303.05
$ws.Range("H7").Value = 347.37
$ws.Range("J7").Value = 422.9
$ws.Range("X7").Value = 312.89
$ws.Range("Y7").Value = 302.93
$ws.Range("Z7").Value = 367.25
$ws.Range("AA7").Value = 384.43
$ws.Range("AB7").Value = 368.34
$ws.Range("AJ7").Value = 423.14
$ws.Range("AU7").Value = 403.78
$ws.Range("AW7").Value = 442.11
# Remove the now-unused trailing blank rows (10-15)
$ws.Range("A10:BA15").Delete()

Write-Output "done"
